$d = $word.ActiveDocument

# --- 1. Add strikethrough formatting to "Actualizar contenido página principal." ---
$rngStrike = $d.Content
$null = $rngStrike.Find.Execute("Actualizar contenido página principal.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngStrike.Font.StrikeThrough = 1

# --- 2. Move the "_GoBack" bookmark so it wraps
#        "Actualizar contenido página principal. " (run + following space run) ---
$rngBookmark = $d.Content
$null = $rngBookmark.Find.Execute("Actualizar contenido página principal. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $rngBookmark)
